# Updated cryptos list on Thu Aug  3 01:19:31 UTC 2023 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) figures for the crypto table on Sheet1.
# Numeric-looking price strings are written through a Text-formatted cell
# (NumberFormat "@" + Value2) so they stay literal text (e.g. "160.20",
# "1.000") instead of being coerced into floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.175.99'
$ws.Range('E2').Value = '  -2.26%  '

$ws.Range('D3').Value = '1.838.58'
$ws.Range('E3').Value = '  -1.72%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value2 = '0.9990'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value2 = '240.04'
$ws.Range('E5').Value = '  -2.70%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value2 = '0.6830'
$ws.Range('E6').Value = '  -2.72%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value2 = '0.9997'
$ws.Range('E7').Value = '  +0.12%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value2 = '0.2996'
$ws.Range('E8').Value = '  -3.33%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value2 = '0.07448'
$ws.Range('E9').Value = '  -4.31%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value2 = '23.17'
$ws.Range('E10').Value = '  -3.72%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value2 = '0.07638'
$ws.Range('E11').Value = '  -2.68%  '

$ws.Range('D12').Value = '1.835.87'
$ws.Range('E12').Value = '  -1.83%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value2 = '5.047'
$ws.Range('E13').Value = '  -2.81%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value2 = '0.6801'
$ws.Range('E14').Value = '  -2.57%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value2 = '87.66'
$ws.Range('E15').Value = '  -6.04%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value2 = '6.136'
$ws.Range('E16').Value = '  -7.76%  '

$ws.Range('D17').Value = '29.166.56'
$ws.Range('E17').Value = '  -2.14%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value2 = '0.000008183'
$ws.Range('E18').Value = '  -2.91%  '

$ws.Range('D19').Value = '2.079.25'
$ws.Range('E19').Value = '  -1.59%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value2 = '228.61'
$ws.Range('E20').Value = '  -6.55%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value2 = '12.52'
$ws.Range('E21').Value = '  -2.89%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value2 = '1.000'
$ws.Range('E22').Value = '  +0.10%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value2 = '7.344'
$ws.Range('E23').Value = '  -3.75%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value2 = '0.9994'
$ws.Range('E24').Value = '  +0.02%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value2 = '160.20'
$ws.Range('E25').Value = '  -0.04%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value2 = '0.1447'
$ws.Range('E26').Value = '  -5.10%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value2 = '8.701'
$ws.Range('E27').Value = '  -3.11%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value2 = '18.04'
$ws.Range('E28').Value = '  -2.42%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value2 = '1.505'
$ws.Range('E29').Value = '  -2.81%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value2 = '4.286'
$ws.Range('E30').Value = '  +0.04%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value2 = '4.147'
$ws.Range('E31').Value = '  -2.66%  '

$ws.Range('E32').Value = '  -0.90%  '

$ws.Range('E33').Value = '  +1.84%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value2 = '0.7541'
$ws.Range('E34').Value = '  -4.86%  '

$ws.Range('E35').Value = '  -4.44%  '

$ws.Range('E36').Value = '  -3.01%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value2 = '2.689'
$ws.Range('E37').Value = '  -0.50%  '

$ws.Range('D38').Value = '1.304.47'
$ws.Range('E38').Value = '  -2.80%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value2 = '0.01827'
$ws.Range('E39').Value = '  -3.35%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value2 = '2.716'
$ws.Range('E40').Value = '  -1.13%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value2 = '0.9498'
$ws.Range('E41').Value = '  -0.96%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value2 = '6.016'
$ws.Range('E42').Value = '  -1.48%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value2 = '104.74'
$ws.Range('E43').Value = '  -2.66%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value2 = '0.9993'
$ws.Range('E44').Value = '  +0.09%  '

# Rows 45-47 reshuffled: RocketPoolETH / Mantle / BabyDogeCoin swap order
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.974.99'
$ws.Range('E45').Value = '  -1.84%  '

$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value2 = '0.5190'
$ws.Range('E46').Value = '  -0.14%  '

$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value2 = '0.00000000122'
$ws.Range('E47').Value = '  -2.54%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value2 = '64.45'
$ws.Range('E48').Value = '  -2.62%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value2 = '9.455'
$ws.Range('E49').Value = '  -4.26%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value2 = '1.768'
$ws.Range('E50').Value = '  -1.45%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value2 = '0.07568'
$ws.Range('E51').Value = '  +16.25%  '
